# Apply corrected temperature-derived mole-fraction / uncertainty values
# to the Ethanol_AVE_MOL_FRAC workbook (Mole_Fractions + Uncertainties sheets),
# and widen the affected numeric columns on Mole_Fractions to fit the new values.

$wb = $excel.ActiveWorkbook

$wsMoleFrac = $wb.Worksheets.Item("Mole_Fractions")
$wsMoleFrac.Range("B2").Value = 0.050975929752124817
$wsMoleFrac.Range("C2").Value = 0.032894614408185391
$wsMoleFrac.Range("D2").Value = 0.16836536701845861
$wsMoleFrac.Range("E2").Value = 0.0054161568256995253
$wsMoleFrac.Range("F2").Value = 0.019995854636702816
$wsMoleFrac.Range("I2").Value = 0.47997688043695741
$wsMoleFrac.Range("J2").Value = 0.14698142924619895
$wsMoleFrac.Range("K2").Value = 0.011128127603664975
$wsMoleFrac.Range("M2").Value = 0.053669451241893054
$wsMoleFrac.Range("N2").Value = 0.018786518943293605
$wsMoleFrac.Range("P2").Value = 0.00022985713744390121
$wsMoleFrac.Range("Q2").Value = 0.00046181892544300859
$wsMoleFrac.Range("R2").Value = 0.0071767355593370307
$wsMoleFrac.Range("B3").Value = 0.054221547116626391
$wsMoleFrac.Range("C3").Value = 0.037792591831608539
$wsMoleFrac.Range("D3").Value = 0.17948627359068156
$wsMoleFrac.Range("E3").Value = 0.0060500332791019419
$wsMoleFrac.Range("F3").Value = 0.022728370721203074
$wsMoleFrac.Range("I3").Value = 0.54033888385867945
$wsMoleFrac.Range("J3").Value = 0.056514356592965199
$wsMoleFrac.Range("K3").Value = 0.015514099177698366
$wsMoleFrac.Range("M3").Value = 0.057783957566546348
$wsMoleFrac.Range("N3").Value = 0.017307771465087711
$wsMoleFrac.Range("P3").Value = 0.00024725149765762418
$wsMoleFrac.Range("Q3").Value = 0.00038506705611319857
$wsMoleFrac.Range("R3").Value = 0.0071989983618921801
$wsMoleFrac.Range("B4").Value = 0.057997323405695682
$wsMoleFrac.Range("C4").Value = 0.033797078629427865
$wsMoleFrac.Range("D4").Value = 0.16781254443338986
$wsMoleFrac.Range("E4").Value = 0.0065530329628752642
$wsMoleFrac.Range("F4").Value = 0.026454497454176606
$wsMoleFrac.Range("I4").Value = 0.58290383808986768
$wsMoleFrac.Range("J4").Value = 0.025371649289658585
$wsMoleFrac.Range("K4").Value = 0.015151198761170668
$wsMoleFrac.Range("M4").Value = 0.057012489818122604
$wsMoleFrac.Range("N4").Value = 0.015985380046762032
$wsMoleFrac.Range("P4").Value = 0.00023016058927327463
$wsMoleFrac.Range("Q4").Value = 0.0003351338336359471
$wsMoleFrac.Range("R4").Value = 0.0065141913376529931
$wsMoleFrac.Range("B5").Value = 0.056191491988095428
$wsMoleFrac.Range("C5").Value = 0.023217706483662559
$wsMoleFrac.Range("D5").Value = 0.15595377452024017
$wsMoleFrac.Range("E5").Value = 0.0072280441159209945
$wsMoleFrac.Range("F5").Value = 0.043109687218418409
$wsMoleFrac.Range("I5").Value = 0.6355516999215518
$wsMoleFrac.Range("J5").Value = 0.0064988602843510474
$wsMoleFrac.Range("K5").Value = 0.010421645060556088
$wsMoleFrac.Range("M5").Value = 0.04350373247504128
$wsMoleFrac.Range("N5").Value = 0.0097740000442874571
$wsMoleFrac.Range("P5").Value = 0.00018529850574693533
$wsMoleFrac.Range("Q5").Value = 0.00021372062982712869
$wsMoleFrac.Range("R5").Value = 0.0042618723039325734
$wsMoleFrac.Range("B6").Value = 0.049475471652917963
$wsMoleFrac.Range("C6").Value = 0.019418482368356107
$wsMoleFrac.Range("D6").Value = 0.1239522210178203
$wsMoleFrac.Range("E6").Value = 0.0076581918795310781
$wsMoleFrac.Range("F6").Value = 0.061374317564195119
$wsMoleFrac.Range("I6").Value = 0.67807117990069898
$wsMoleFrac.Range("J6").Value = 0.0034247664756096884
$wsMoleFrac.Range("K6").Value = 0.007836996942937496
$wsMoleFrac.Range("M6").Value = 0.035264053336778914
$wsMoleFrac.Range("N6").Value = 0.0061460220381650719
$wsMoleFrac.Range("P6").Value = 0.000085398442940802849
$wsMoleFrac.Range("Q6").Value = 0.00016109953882472187
$wsMoleFrac.Range("R6").Value = 0.0029780072530081634
$wsMoleFrac.Range("B7").Value = 0.042846117388756745
$wsMoleFrac.Range("C7").Value = 0.008048043593601905
$wsMoleFrac.Range("D7").Value = 0.086275401159931195
$wsMoleFrac.Range("E7").Value = 0.0082133802543575621
$wsMoleFrac.Range("F7").Value = 0.09722745784382697
$wsMoleFrac.Range("I7").Value = 0.728032039764831
$wsMoleFrac.Range("J7").Value = 0.00068403962433076005
$wsMoleFrac.Range("K7").Value = 0.0033480435359149398
$wsMoleFrac.Range("M7").Value = 0.01721106512990201
$wsMoleFrac.Range("N7").Value = 0.0027306867369380219
$wsMoleFrac.Range("Q7").Value = 0.000079356911468174106
$wsMoleFrac.Range("R7").Value = 0.0011624914740394965
$wsMoleFrac.Range("B8").Value = 0.034860829869311966
$wsMoleFrac.Range("C8").Value = 0.0024038317817637353
$wsMoleFrac.Range("D8").Value = 0.061624935810477427
$wsMoleFrac.Range("E8").Value = 0.0084699707007492506
$wsMoleFrac.Range("F8").Value = 0.12696471988061858
$wsMoleFrac.Range("I8").Value = 0.75351305989515094
$wsMoleFrac.Range("J8").Value = 0.000081207887185947791
$wsMoleFrac.Range("K8").Value = 0.00091278296012955833
$wsMoleFrac.Range("M8").Value = 0.0060808910970214331
$wsMoleFrac.Range("N8").Value = 0.00086886923762618726
$wsMoleFrac.Range("Q8").Value = 0.00001332742460308236
$wsMoleFrac.Range("R8").Value = 0.00026681163584672493
$wsMoleFrac.Range("B9").Value = 0.021730777368541822
$wsMoleFrac.Range("D9").Value = 0.029775700234771479
$wsMoleFrac.Range("E9").Value = 0.0086100342900754315
$wsMoleFrac.Range("F9").Value = 0.16151466931957065
$wsMoleFrac.Range("I9").Value = 0.773840729790777
$wsMoleFrac.Range("J9").Value = 0.0000056225724516888431
$wsMoleFrac.Range("K9").Value = 0.000038937102597632847
$wsMoleFrac.Range("M9").Value = 0.00034930317882126902
$wsMoleFrac.Range("R9").Value = 0.00010513634682183687
$wsMoleFrac.Range("B10").Value = 0.01315828095408597
$wsMoleFrac.Range("D10").Value = 0.016377860457377297
$wsMoleFrac.Range("E10").Value = 0.0087684504772820668
$wsMoleFrac.Range("F10").Value = 0.18275741027607564
$wsMoleFrac.Range("I10").Value = 0.7745488754110702
$wsMoleFrac.Range("K10").Value = 0.000018362471229638906
$wsMoleFrac.Range("M10").Value = 0.000020214073193780772

$wsUncert = $wb.Worksheets.Item("Uncertainties")
$wsUncert.Range("B2").Value = 0.0026266193630667583
$wsUncert.Range("C2").Value = 0.0021672178900110375
$wsUncert.Range("D2").Value = 0.033824572086532077
$wsUncert.Range("E2").Value = 0.00066923849323026095
$wsUncert.Range("F2").Value = 0.0074280214850025745
$wsUncert.Range("I2").Value = 0.046322573535070126
$wsUncert.Range("J2").Value = 0.06987942925734722
$wsUncert.Range("K2").Value = 0.00077332619957453134
$wsUncert.Range("M2").Value = 0.0037869151164321648
$wsUncert.Range("N2").Value = 0.0034057587376716621
$wsUncert.Range("P2").Value = 0.00010599816278823837
$wsUncert.Range("Q2").Value = 0.000035872330925911651
$wsUncert.Range("R2").Value = 0.0010190164167024254
$wsUncert.Range("B3").Value = 0.0057055855478190583
$wsUncert.Range("C3").Value = 0.0075288976184398884
$wsUncert.Range("D3").Value = 0.051343038014147899
$wsUncert.Range("E3").Value = 0.00067512478297217964
$wsUncert.Range("F3").Value = 0.018150556176502235
$wsUncert.Range("I3").Value = 0.050150779184922882
$wsUncert.Range("J3").Value = 0.022962732202062401
$wsUncert.Range("K3").Value = 0.0044387627801717169
$wsUncert.Range("M3").Value = 0.012365153713955945
$wsUncert.Range("N3").Value = 0.0046906711199488629
$wsUncert.Range("P3").Value = 0.000054853820526955779
$wsUncert.Range("Q3").Value = 0.00011890309927481706
$wsUncert.Range("R3").Value = 0.0019789773874642458
$wsUncert.Range("B4").Value = 0.0036773902397962068
$wsUncert.Range("C4").Value = 0.0025161065915431792
$wsUncert.Range("D4").Value = 0.05406749018361566
$wsUncert.Range("E4").Value = 0.00059716560115085097
$wsUncert.Range("F4").Value = 0.0017828360570605042
$wsUncert.Range("I4").Value = 0.03618940880327752
$wsUncert.Range("J4").Value = 0.0050620291944718015
$wsUncert.Range("K4").Value = 0.0016751729630561738
$wsUncert.Range("M4").Value = 0.0035583663309037036
$wsUncert.Range("N4").Value = 0.0017080450018450519
$wsUncert.Range("P4").Value = 0.000065902821011222854
$wsUncert.Range("Q4").Value = 0.00003012937664850317
$wsUncert.Range("R4").Value = 0.00052814925313633145
$wsUncert.Range("B5").Value = 0.0027112364473577876
$wsUncert.Range("C5").Value = 0.0059245851097677292
$wsUncert.Range("D5").Value = 0.031057346856867254
$wsUncert.Range("E5").Value = 0.00070450613975044118
$wsUncert.Range("F5").Value = 0.013839846225593923
$wsUncert.Range("I5").Value = 0.046049782052465765
$wsUncert.Range("J5").Value = 0.0014908636566009629
$wsUncert.Range("K5").Value = 0.0018247107196026507
$wsUncert.Range("M5").Value = 0.0074436747135353748
$wsUncert.Range("N5").Value = 0.00070073542577894375
$wsUncert.Range("P5").Value = 0.000051173530448865387
$wsUncert.Range("Q5").Value = 0.000036355222491489642
$wsUncert.Range("R5").Value = 0.00081320519115205307
$wsUncert.Range("B6").Value = 0.0077476534833186049
$wsUncert.Range("C6").Value = 0.015416524965127635
$wsUncert.Range("D6").Value = 0.028313838985886217
$wsUncert.Range("E6").Value = 0.00070219285578005438
$wsUncert.Range("F6").Value = 0.041376483343277086
$wsUncert.Range("I6").Value = 0.044760494554843669
$wsUncert.Range("J6").Value = 0.0043329053048977223
$wsUncert.Range("K6").Value = 0.0060026829054579444
$wsUncert.Range("M6").Value = 0.02367082852072009
$wsUncert.Range("N6").Value = 0.0034384257609070164
$wsUncert.Range("P6").Value = 0.000052197355625829127
$wsUncert.Range("Q6").Value = 0.000030562589856713837
$wsUncert.Range("R6").Value = 0.0024205394934442175
$wsUncert.Range("B7").Value = 0.0011799693942510226
$wsUncert.Range("C7").Value = 0.0003986895141974132
$wsUncert.Range("D7").Value = 0.017199794792909086
$wsUncert.Range("E7").Value = 0.00059216814337112434
$wsUncert.Range("F7").Value = 0.0024477695663818839
$wsUncert.Range("I7").Value = 0.018231784058131924
$wsUncert.Range("J7").Value = 0.000061665550613785415
$wsUncert.Range("K7").Value = 0.00018193794934851098
$wsUncert.Range("M7").Value = 0.00046715487069406433
$wsUncert.Range("N7").Value = 0.000055323291367082532
$wsUncert.Range("Q7").Value = 0.000003646859755837121
$wsUncert.Range("R7").Value = 0.000065191065175053078
$wsUncert.Range("B8").Value = 0.0015914950694957896
$wsUncert.Range("C8").Value = 0.000599279707907785
$wsUncert.Range("D8").Value = 0.019829253276285805
$wsUncert.Range("E8").Value = 0.00069949377489719234
$wsUncert.Range("F8").Value = 0.0099227846194020943
$wsUncert.Range("I8").Value = 0.037603404511584893
$wsUncert.Range("J8").Value = 0.00006095641957567591
$wsUncert.Range("K8").Value = 0.00028607045751298906
$wsUncert.Range("M8").Value = 0.0015419912670630008
$wsUncert.Range("N8").Value = 0.00028454784765668693
$wsUncert.Range("Q8").Value = 0.000026671960725986185
$wsUncert.Range("R8").Value = 0.000095214341229518979
$wsUncert.Range("B9").Value = 0.0021981550126556038
$wsUncert.Range("D9").Value = 0.0087908298441017902
$wsUncert.Range("E9").Value = 0.00063598178520216737
$wsUncert.Range("F9").Value = 0.0080419733877810547
$wsUncert.Range("I9").Value = 0.025348595333920886
$wsUncert.Range("J9").Value = 0.0000040675972841788058
$wsUncert.Range("K9").Value = 0.0000086642898628642752
$wsUncert.Range("M9").Value = 0.00069869783803481599
$wsUncert.Range("R9").Value = 0.00021042676296313123
$wsUncert.Range("B10").Value = 0.0012304227947352738
$wsUncert.Range("D10").Value = 0.0045582212887324983
$wsUncert.Range("E10").Value = 0.00071263271454298548
$wsUncert.Range("F10").Value = 0.0086858388377742003
$wsUncert.Range("I10").Value = 0.024920378562480826
$wsUncert.Range("K10").Value = 0.000036750900487032914
$wsUncert.Range("M10").Value = 0.000040434331527054323

# Columns B:F, I, K, M on Mole_Fractions now hold longer decimal literals;
# widen them to 12 characters (matches the workbook author's resize).
$wsMoleFrac.Columns.Item(2).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(3).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(4).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(5).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(6).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(9).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(11).ColumnWidth = 11.14
$wsMoleFrac.Columns.Item(13).ColumnWidth = 11.14
